$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The RANK.AVG formulas in column E (rows 2-7) were computed against the
# whole C2:C13 range before the B (coverage) column existed; those rows
# don't participate in the corrected correlation, so their ranks are removed.
$ws.Range("E2").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("E7").ClearContents()

# Fill in the "Overall Statement Coverage" values for the rows that are
# actually used in the Spearman correlation (rows 8-12); row 13's old value
# (0.89) belongs to row 12 now, so row13 (B13) is left blank.
$ws.Range("B8").Value = 0.85
$ws.Range("B9").Value = 0.84
$ws.Range("B10").Value = 0.86
$ws.Range("B11").Value = 0.86
$ws.Range("B12").Value = 0.89
$ws.Range("B13").ClearContents()

# Rank of statement coverage (column D) for rows 8-12.
$ws.Range("D8").Formula = "=RANK.AVG(B8,`$B`$8:`$B`$13,1)"
$ws.Range("D9:D12").Formula = "=RANK.AVG(B9,`$B`$8:`$B`$13,1)"

# Corrected rank of code churn (column E) — ranked only over C8:C12 now.
$ws.Range("E8").Formula = "=RANK.AVG(C8,`$C`$8:`$C`$12,1)"
$ws.Range("E9:E12").Formula = "=RANK.AVG(C9,`$C`$8:`$C`$12,1)"

# Square of the rank difference (column F).
$ws.Range("F8").Formula = "=(E8-D8)^2"
$ws.Range("F9:F12").Formula = "=(E9-D9)^2"

# Sum of squared differences + Spearman's rank correlation coefficient.
$ws.Range("E13").Value = "Sum:"
$ws.Range("F13").Formula = "=SUM(F8:F12)"

$ws.Range("E14").Value = "Spearman:"
$ws.Range("F14").Formula = "=1-((6*F13)/(125-5))"

# Match the saved selection/cursor position.
$ws.Range("F15").Select()
